$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "transformation.py" block (rows 20-23) gains two new functions:
# first_order_derivative and second_order_derivative, inserted right after
# "snv" (row 23). Copy row 23 and use "Insert Copied Cells" so the new rows
# inherit its column-E formatting (style index 1), then overwrite the values.
$ws.Rows("23:23").Copy()
$ws.Rows("24:24").Insert()
$ws.Rows("23:23").Copy()
$ws.Rows("25:25").Insert()

$ws.Range("E24").Value = "first_order_derivative"
$ws.Range("F24").Value = "f"

$ws.Range("E25").Value = "second_order_derivative"
$ws.Range("F25").Value = "f"

$excel.CutCopyMode = $false

# Reflect the saved view state: scrolled so row 4 is at the top, with F29
# selected as the active cell.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("F29").Select()
